# B6-PowerPoint.pptx edit — Sun, Jul 12, 2020  3:05:14 AM
#
# 1) Three tables (on the slides holding the "component 3" grids) get their
#    table style switched from the old GUID to the new one.
# 2) The deck's theme ("Integral" master theme, colour scheme "Red Violet")
#    is swapped for the "Office Theme" colour palette (colour scheme
#    "Office") that used to live only on the notes master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables (slides 14, 15, 16 each hold one table as
#    their first shape) with the new table-style id.
# ---------------------------------------------------------------------
$newTableStyleId = "{2E00D5F9-5963-4D7E-A0C9-BAE4AA9A8B87}"

for ($i = 14; $i -le 16; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId, $true)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colours used by the slide master from the "Red
#    Violet" Integral palette to the standard "Office" palette. PowerPoint
#    stores RGB as 0x00BBGGRR, so convert each target hex colour (RRGGBB)
#    accordingly before assigning it.
# ---------------------------------------------------------------------
function HexToBgr($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Order matches ThemeColorScheme.Colors(1..12): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $targetHex = $officeColors[$i - 1]
    $themeColors.Colors($i).RGB = HexToBgr $targetHex
}
